# Refresh the Siglec1-Spn ligand-receptor table with newly recomputed TPM
# values. The sending/target cluster set changed too: "ECs" is a new
# sending cluster and "MuSCs" moved from being a target cluster to a
# sending cluster (it no longer appears as a target). Columns are:
# A Sending cluster, B Ligand symbol, C Receptor symbol, D Target cluster,
# E..J ligand stats, K..P receptor stats, Q..T edge weights/specificities.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
  @("ECs", "Siglec1", "Spn", "Inflammatory-Mac", 2, 0.6666666666666666, 0.5705793333333333, 1.711738, 0.00495357279640748, 0.004953572796407481, 2, 0.6666666666666666, 3.260090333333334, 9.780271000000001, 0.5507055085354173, 0.5507055085354173, 1.860140168999778, 16.741261520998, 0.00272795982591279, 0.002727959825912791),
  @("ECs", "Siglec1", "Spn", "Neutrophils", 2, 0.6666666666666666, 0.5705793333333333, 1.711738, 0.00495357279640748, 0.004953572796407481, 2, 0.6666666666666666, 1.435229, 4.305687, 0.2424437471036677, 0.2424437471036677, 0.8189120060006666, 7.370208054006, 0.001200962750311823, 0.001200962750311823),
  @("ECs", "Siglec1", "Spn", "Resolving-Mac", 2, 0.6666666666666666, 0.5705793333333333, 1.711738, 0.00495357279640748, 0.004953572796407481, 3, 1, 1.224524, 3.673572, 0.2068507443609149, 0.2068507443609149, 0.6986880875706666, 6.288192788136, 0.001024650220182866, 0.001024650220182866),
  @("FAPs", "Siglec1", "Spn", "Inflammatory-Mac", 2, 0.6666666666666666, 0.04769933333333334, 0.143098, 0.0004141091452198395, 0.0004141091452198395, 2, 0.6666666666666666, 3.260090333333334, 9.780271000000001, 0.5507055085354173, 0.5507055085354173, 0.1555041355064445, 1.399537219558, 0.0002280521874074587, 0.0002280521874074587),
  @("FAPs", "Siglec1", "Spn", "Neutrophils", 2, 0.6666666666666666, 0.04769933333333334, 0.143098, 0.0004141091452198395, 0.0004141091452198395, 2, 0.6666666666666666, 1.435229, 4.305687, 0.2424437471036677, 0.2424437471036677, 0.06845946648066667, 0.6161351983259999, 0.0001003981728769948, 0.0001003981728769948),
  @("FAPs", "Siglec1", "Spn", "Resolving-Mac", 2, 0.6666666666666666, 0.04769933333333334, 0.143098, 0.0004141091452198395, 0.0004141091452198395, 3, 1, 1.224524, 3.673572, 0.2068507443609149, 0.2068507443609149, 0.05840897845066667, 0.525680806056, 0.00008565878493538602, 0.00008565878493538603),
  @("Inflammatory-Mac", "Siglec1", "Spn", "Inflammatory-Mac", 3, 1, 37.411254, 112.233762, 0.3247915921021042, 0.3247915921021042, 2, 0.6666666666666666, 3.260090333333334, 9.780271000000001, 0.5507055085354173, 0.5507055085354173, 121.964067523278, 1097.676607709502, 0.1788645188966171, 0.1788645188966171),
  @("Inflammatory-Mac", "Siglec1", "Spn", "Neutrophils", 3, 1, 37.411254, 112.233762, 0.3247915921021042, 0.3247915921021042, 2, 0.6666666666666666, 1.435229, 4.305687, 0.2424437471036677, 0.2424437471036677, 53.69371666716599, 483.243450004494, 0.07874369061700014, 0.07874369061700014),
  @("Inflammatory-Mac", "Siglec1", "Spn", "Resolving-Mac", 3, 1, 37.411254, 112.233762, 0.3247915921021042, 0.3247915921021042, 3, 1, 1.224524, 3.673572, 0.2068507443609149, 0.2068507443609149, 45.810978393096, 412.298805537864, 0.06718338258848691, 0.06718338258848691),
  @("MuSCs", "Siglec1", "Spn", "Inflammatory-Mac", 1, 0.3333333333333333, 0.021226, 0.063678, 0.0001842768043530234, 0.0001842768043530234, 2, 0.6666666666666666, 3.260090333333334, 9.780271000000001, 0.5507055085354173, 0.5507055085354173, 0.06919867741533334, 0.6227880967380001, 0.0001014822512525133, 0.0001014822512525133),
  @("MuSCs", "Siglec1", "Spn", "Neutrophils", 1, 0.3333333333333333, 0.021226, 0.063678, 0.0001842768043530234, 0.0001842768043530234, 2, 0.6666666666666666, 1.435229, 4.305687, 0.2424437471036677, 0.2424437471036677, 0.03046417075399999, 0.274177536786, 0.00004467675895163644, 0.00004467675895163645),
  @("MuSCs", "Siglec1", "Spn", "Resolving-Mac", 1, 0.3333333333333333, 0.021226, 0.063678, 0.0001842768043530234, 0.0001842768043530234, 3, 1, 1.224524, 3.673572, 0.2068507443609149, 0.2068507443609149, 0.025991746424, 0.233925717816, 0.00003811779414887357, 0.00003811779414887358),
  @("Neutrophils", "Siglec1", "Spn", "Inflammatory-Mac", 3, 1, 12.49536033333333, 37.486081, 0.1084804047614339, 0.1084804047614339, 2, 0.6666666666666666, 3.260090333333334, 9.780271000000001, 0.5507055085354173, 0.5507055085354173, 40.73600343421678, 366.624030907951, 0.05974075647027334, 0.05974075647027335),
  @("Neutrophils", "Siglec1", "Spn", "Neutrophils", 3, 1, 12.49536033333333, 37.486081, 0.1084804047614339, 0.1084804047614339, 2, 0.6666666666666666, 1.435229, 4.305687, 0.2424437471036677, 0.2424437471036677, 17.93370351584966, 161.403331642647, 0.02630039581768458, 0.02630039581768458),
  @("Neutrophils", "Siglec1", "Spn", "Resolving-Mac", 3, 1, 12.49536033333333, 37.486081, 0.1084804047614339, 0.1084804047614339, 3, 1, 1.224524, 3.673572, 0.2068507443609149, 0.2068507443609149, 15.30086861681466, 137.707817551332, 0.02243925247347594, 0.02243925247347594),
  @("Resolving-Mac", "Siglec1", "Spn", "Inflammatory-Mac", 3, 1, 64.639295, 193.917885, 0.5611760443904816, 0.5611760443904816, 2, 0.6666666666666666, 3.260090333333334, 9.780271000000001, 0.5507055085354173, 0.5507055085354173, 210.7299407829817, 1896.569467046835, 0.3090427389039541, 0.3090427389039541),
  @("Resolving-Mac", "Siglec1", "Spn", "Neutrophils", 3, 1, 64.639295, 193.917885, 0.5611760443904816, 0.5611760443904816, 2, 0.6666666666666666, 1.435229, 4.305687, 0.2424437471036677, 0.2424437471036677, 92.772190723555, 834.949716511995, 0.1360536229868425, 0.1360536229868425),
  @("Resolving-Mac", "Siglec1", "Spn", "Resolving-Mac", 3, 1, 64.639295, 193.917885, 0.5611760443904816, 0.5611760443904816, 3, 1, 1.224524, 3.673572, 0.2068507443609149, 0.2068507443609149, 79.15236807058, 712.3713126352201, 0.116079682499685, 0.116079682499685)
)

$r = 2
foreach ($row in $data) {
  $c = 1
  foreach ($val in $row) {
    $ws.Cells.Item($r, $c).Value = $val
    $c = $c + 1
  }
  $r = $r + 1
}
